$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Forecast Comparison")
$ws2 = $wb.Worksheets.Item("Summary")

# --- Forecast Comparison sheet: updated forecast figures (Prophet / Amazon Mean / P70 / P80 / P90) ---
$ws1.Range("C2").Value = 918
$ws1.Range("D2").Value = 842
$ws1.Range("E2").Value = 1011
$ws1.Range("F2").Value = 1189
$ws1.Range("G2").Value = 1467
$ws1.Range("C3").Value = 932
$ws1.Range("D3").Value = 826
$ws1.Range("E3").Value = 1000
$ws1.Range("F3").Value = 1195
$ws1.Range("G3").Value = 1504
$ws1.Range("C4").Value = 844
$ws1.Range("D4").Value = 672
$ws1.Range("E4").Value = 813
$ws1.Range("F4").Value = 968
$ws1.Range("G4").Value = 1215
$ws1.Range("C5").Value = 712
$ws1.Range("D5").Value = 555
$ws1.Range("E5").Value = 669
$ws1.Range("F5").Value = 793
$ws1.Range("G5").Value = 988
$ws1.Range("C6").Value = 612
$ws1.Range("D6").Value = 565
$ws1.Range("E6").Value = 683
$ws1.Range("F6").Value = 816
$ws1.Range("G6").Value = 1025
$ws1.Range("C7").Value = 581
$ws1.Range("D7").Value = 542
$ws1.Range("E7").Value = 657
$ws1.Range("F7").Value = 785
$ws1.Range("G7").Value = 989
$ws1.Range("C8").Value = 602
$ws1.Range("D8").Value = 511
$ws1.Range("E8").Value = 622
$ws1.Range("F8").Value = 754
$ws1.Range("G8").Value = 966
$ws1.Range("C9").Value = 647
$ws1.Range("D9").Value = 513
$ws1.Range("E9").Value = 622
$ws1.Range("F9").Value = 747
$ws1.Range("G9").Value = 944
$ws1.Range("C10").Value = 688
$ws1.Range("D10").Value = 480
$ws1.Range("E10").Value = 581
$ws1.Range("F10").Value = 695
$ws1.Range("G10").Value = 875
$ws1.Range("C11").Value = 704
$ws1.Range("D11").Value = 464
$ws1.Range("E11").Value = 562
$ws1.Range("F11").Value = 671
$ws1.Range("G11").Value = 845
$ws1.Range("C12").Value = 676
$ws1.Range("D12").Value = 453
$ws1.Range("E12").Value = 552
$ws1.Range("F12").Value = 673
$ws1.Range("G12").Value = 867
$ws1.Range("C13").Value = 589
$ws1.Range("D13").Value = 440
$ws1.Range("E13").Value = 537
$ws1.Range("F13").Value = 660
$ws1.Range("G13").Value = 859
$ws1.Range("C14").Value = 458
$ws1.Range("D14").Value = 482
$ws1.Range("E14").Value = 588
$ws1.Range("F14").Value = 716
$ws1.Range("G14").Value = 921
$ws1.Range("C15").Value = 343
$ws1.Range("D15").Value = 454
$ws1.Range("E15").Value = 554
$ws1.Range("F15").Value = 680
$ws1.Range("G15").Value = 884
$ws1.Range("C16").Value = 285
$ws1.Range("D16").Value = 450
$ws1.Range("E16").Value = 550
$ws1.Range("F16").Value = 677
$ws1.Range("G16").Value = 881
$ws1.Range("C17").Value = 284
$ws1.Range("D17").Value = 432
$ws1.Range("E17").Value = 527
$ws1.Range("F17").Value = 653
$ws1.Range("G17").Value = 858

# --- Summary sheet: recompute dependent summary metrics (kept as text to match existing cell formatting) ---
function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws2.Range("B9") "9875"
Set-TextValue $ws2.Range("B10") "5848"
Set-TextValue $ws2.Range("B11") "3406"
Set-TextValue $ws2.Range("B12") "932"
Set-TextValue $ws2.Range("B14") "284"
Set-TextValue $ws2.Range("B15") "2025-03-23"
